$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin prices (column D) and 1h volume deltas (column E) scraped by the
# GitHub Actions cron job.
#
# A handful of the new D-column prices are plain decimal text (e.g. "4.36",
# "147.30", "0.530") that Excel would otherwise silently reinterpret as a
# floating-point number (losing the exact trailing-zero text and diverging from
# the source). For those cells we briefly switch the cell to a Text number
# format before assigning the value, then restore the cell style to Normal so
# no lasting formatting change is left behind.

$ws.Range("D2").Value = "26.850.56"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.642.32"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.495"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.93%  "
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("E9").Value = "  -0.96%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("D12").Value = "1.871.94"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("D13").Value = "1.640.65"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.20%  "
$ws.Range("D17").Value = "26.851.14"
$ws.Range("D18").Value = "0.0₃0730"
$ws.Range("E18").Value = "  -0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "215.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.08%  "
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0508"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("E31").Value = "  +0.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.41%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("D34").Value = "1.281.42"
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("E35").Value = "  +1.33%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.530"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.820"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.22%  "
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("D43").Value = "1.782.75"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.09"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.59%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "61.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0968"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.53%  "
$ws.Range("E51").Value = "  +0.01%  "
